# Fix Import-Function: rename existing Region/Level columns (A/E) to the
# "_Chiara" variants, and add three freshly (re-)computed columns F:H
# ("add" flag, numeric "Level", text "Region") used for the CAA talk.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing header cells ---
$ws.Range("A1").Value = "Region_Chiara"
$ws.Range("E1").Value = "Level_Chiara"

# --- Headers for the three newly appended columns ---
$ws.Range("F1").Value = "add"
$ws.Range("G1").Value = "Level"
$ws.Range("H1").Value = "Region"

# Plain Range.Value assignment of the strings "TRUE"/"FALSE" gets auto-coerced
# to a Boolean cell by Excel, but column F must hold literal text. Build one
# reusable text "TRUE" cell and one text "FALSE" cell off to the side (via a
# TEXT() formula + paste-values), then stamp copies of those into column F.
$helper = $ws.Range("Z1")
$helper.Formula = '=TEXT(TRUE,"@")'
$helper.Copy()
$ws.Range("Z2").PasteSpecial(-4163)
$helper.Formula = '=TEXT(FALSE,"@")'
$helper.Copy()
$ws.Range("Z3").PasteSpecial(-4163)
$textTrue = $ws.Range("Z2")
$textFalse = $ws.Range("Z3")

# --- Per-row data for the new F:H columns ---
$rows = @(
    @{ r = 2; f = "FALSE"; g = 1; h = "low socio- economic status" },
    @{ r = 3; f = "TRUE"; g = 1; h = "low socio- economic status" },
    @{ r = 4; f = "FALSE"; g = 1; h = "high socio- economic status" },
    @{ r = 5; f = "FALSE"; g = 2; h = "high socio- economic status" },
    @{ r = 6; f = "FALSE"; g = 2; h = "high socio- economic status" },
    @{ r = 7; f = "FALSE"; g = 1; h = "high socio- economic status" },
    @{ r = 8; f = "FALSE"; g = 1; h = "urban" },
    @{ r = 9; f = "FALSE"; g = 2; h = "urban" },
    @{ r = 10; f = "FALSE"; g = 2; h = "urban" },
    @{ r = 11; f = "TRUE"; g = 1; h = "urban" },
    @{ r = 12; f = "FALSE"; g = 1; h = "plague deaths" },
    @{ r = 13; f = "FALSE"; g = 1; h = "plague deaths" },
    @{ r = 14; f = "FALSE"; g = 1; h = "plague deaths" },
    @{ r = 15; f = "FALSE"; g = 1; h = "plague deaths" },
    @{ r = 16; f = "FALSE"; g = 1; h = "plague deaths" },
    @{ r = 17; f = "FALSE"; g = 1; h = "plague deaths" },
    @{ r = 18; f = "FALSE"; g = 1; h = "plague deaths" }
)

foreach ($row in $rows) {
    if ($row.f -eq "TRUE") { $textTrue.Copy() } else { $textFalse.Copy() }
    $ws.Cells.Item($row.r, 6).PasteSpecial(-4163)   # F: add
    $ws.Cells.Item($row.r, 7).Value = $row.g        # G: Level
    $ws.Cells.Item($row.r, 8).Value = $row.h        # H: Region
}

# --- Remove the scratch helper cells ---
$ws.Range("Z1:Z3").Clear()

# --- Match the saved selection state ---
$ws.Range("A7").Select()